# Update cryptos list values (Price and Volume(1h) columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.184.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.69"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.90"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5247"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07699"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.636"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.684.65"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.12"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5631"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8223"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.173.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.669"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.971"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.298"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.523"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05538"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.273"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.481"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.378"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.570"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9536"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.780"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5722"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01598"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.917"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.033.96"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8347"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.796.25"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.059"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4346"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05237"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("E51").Value = "  -4.47%  "
